# Auto-generated Excel COM-interop edit script
# Commit message: "Added space to selenium dsl commands"
#
# Rebuilds Sheet1 and Sheet2 into a 4-row test-case table (TC_001..TC_003 plus
# an unlabeled wikipedia case), rewrites the Selenium DSL command text to add
# the missing spaces (LoadUrl -> Load Url, GetPageTitle -> Get Page Title,
# MaximizeBrowserWindow -> Maximize Browser Window, QuitBrowser -> Quit
# Browser), adds an "Invoke Browser" first step to the admin-login case, and
# drops the stray trailing decorative rows/borders that Sheet2 had.

$wb = $excel.ActiveWorkbook

$xlHAlignLeft = -4131
$xlVAlignTop = -4160
$xlVAlignCenter = -4108

$B0 = "Add 5 and 2 Assign %{result}%`n%{result}% VerifyEqual 7`nPrint %{result}%"
$B2 = "Load Url http://127.0.0.1/wordpress/wp-admin`nMaximize Browser Window`nSleep 2`nGet Page Title Assign %{title}%`n%{title}% VerifyEqual test › Log In`nPrint %{title}%`nQuit Browser"
$B3 = "Invoke Browser`nLoad Url http://127.0.0.1/wordpress/wp-admin`nType admin in element with id = user_login`nType test in element with id = user_pass`nClick element with id = wp-submit`nSleep 2`nGet Page Title Assign %{title}%`n%{title}% VerifyEqual Dashboard ‹ test — WordPress`nQuit Browser"
$B4 = "Load Url http://en.wikipedia.org/wiki/India`nClick element with linkText = View history`nSleep 5`nGet Page Title Assign %{title}%`n%{title}% VerifyEqual India: Revision history - Wikipedia, the free encyclopedia`nQuit Browser"

function Set-TestCaseSheet($ws) {
    # --- Drop the old trailing rows beyond row 4 (decorative / leftover rows) ---
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count + $used.Row - 1
    if ($lastRow -gt 4) {
        $ws.Range("A5:B" + $lastRow).EntireRow.Delete()
    }

    # --- Column widths & default (empty-cell) column formatting ---
    $ws.Columns.Item(1).ColumnWidth = 8.28515625
    $ws.Columns.Item(2).ColumnWidth = 71

    # --- B1 "Add 5 and 2" text is unchanged; fill the new TC_00x labels next ---
    $ws.Range("B1").Value = $B0
    $ws.Range("A1").Value = "TC_001"
    $ws.Range("A2").Value = "TC_002"
    $ws.Range("A3").Value = "TC_003"
    $ws.Range("A4").ClearContents()

    # --- then the long DSL command bodies, in B4, B2, B3 order ---
    $ws.Range("B4").Value = $B4
    $ws.Range("B2").Value = $B2
    $ws.Range("B3").Value = $B3

    # --- Row heights ---
    $ws.Rows.Item(1).RowHeight = 45
    $ws.Rows.Item(2).RowHeight = 105
    $ws.Rows.Item(3).RowHeight = 135
    $ws.Rows.Item(4).RowHeight = 105

    # --- Cell formatting: column A = left/center/wrap, column B = top/wrap ---
    $dataRange = $ws.Range("A1:A4")
    $dataRange.HorizontalAlignment = $xlHAlignLeft
    $dataRange.VerticalAlignment = $xlVAlignCenter
    $dataRange.WrapText = $true

    $textRange = $ws.Range("B1:B4")
    $textRange.VerticalAlignment = $xlVAlignTop
    $textRange.WrapText = $true
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

Set-TestCaseSheet $ws1
Set-TestCaseSheet $ws2

# --- Sheet1 selection: D2 ---
$ws1.Activate()
$ws1.Range("D2").Select()

# --- Sheet2 selection: whole sheet, then re-activate Sheet1 so it stays the active tab ---
$ws2.Cells.Select()
$ws1.Activate()
$ws1.Range("D2").Select()
